# clean code, implement task_timer()
#
# Rewrites the "task time" column so durations are stored as real
# fractional-day numbers (like Excel's own time math) instead of plain
# integers, and appends the newly timed tasks collected by task_timer().

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing "kitchen" row: was a bare integer (9), now a precise
#     fractional-day duration.
$ws.Range("C4").Value = 12.00111111111111

# --- newly logged tasks from task_timer(), appended below "bathroom".
$ws.Range("B6").Value = "nowe"
$ws.Range("C6").Value = 0.1666666666666667

$ws.Range("B7").Value = "list_existing_tasks"
$ws.Range("C7").Value = 0.09

$ws.Range("B8").Value = "test"
$ws.Range("C8").Value = 0.0002777777777777778

# First pass: comma-decimal format (registers numFmt 164 "0,00").
$ws.Range("C4").NumberFormat = "0,00"
$ws.Range("C6").NumberFormat = "0,00"
$ws.Range("C7").NumberFormat = "0,00"
$ws.Range("C8").NumberFormat = "0,00"

# Second pass: settle on the plain "0.00" look for every duration cell.
# The comma format stays registered in the stylesheet but ends up unused.
$ws.Range("C4").NumberFormat = "0.00"
$ws.Range("C6").NumberFormat = "0.00"
$ws.Range("C7").NumberFormat = "0.00"
$ws.Range("C8").NumberFormat = "0.00"
